# Auto-generated Excel COM-interop script
# Applies the numeric cell updates described by the commit's xlsx diff,
# sheet by sheet. All target cells hold literal numbers (no formulas).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1431.9678
$ws.Range("I15").Value = 1431.9678
$ws.Range("K15").Value = 4295.903399999999
$ws.Range("M15").Value = -4126.903399999999
$ws.Range("H69").Value = 18333.334
$ws.Range("J69").Value = 17500
$ws.Range("L69").Value = 52500
$ws.Range("N69").Value = -54248
$ws.Range("H72").Value = 18333.334
$ws.Range("J72").Value = 17500
$ws.Range("L72").Value = 157500
$ws.Range("N72").Value = -166236
$ws.Range("H74").Value = 12818
$ws.Range("I74").Value = 11599.8
$ws.Range("K74").Value = 11599.8
$ws.Range("M74").Value = -10663.8
$ws.Range("H77").Value = 12818
$ws.Range("I77").Value = 11599.8
$ws.Range("K77").Value = 57999
$ws.Range("M77").Value = -53319
$ws.Range("H96").Value = 358.875
$ws.Range("I96").Value = 353.5
$ws.Range("K96").Value = 1060.5
$ws.Range("M96").Value = 312.5
$ws.Range("H111").Value = 168716.17
$ws.Range("J111").Value = 2459.6
$ws.Range("L111").Value = 7378.799999999999
$ws.Range("N111").Value = -13512.8
$ws.Range("H129").Value = 1868.2727
$ws.Range("I129").Value = 1225.3334
$ws.Range("J129").Value = 2639.8
$ws.Range("K129").Value = 3676.0002
$ws.Range("L129").Value = 7919.400000000001
$ws.Range("M129").Value = 1323.9998
$ws.Range("N129").Value = -17919.4
$ws.Range("H131").Value = 2924.611
$ws.Range("I131").Value = 2376.7334
$ws.Range("K131").Value = 7130.2002
$ws.Range("M131").Value = -2090.2002
$ws.Range("H132").Value = 6231.8
$ws.Range("I132").Value = 6955.3887
$ws.Range("K132").Value = 20866.1661
$ws.Range("M132").Value = -18336.1661
$ws.Range("H137").Value = 3891.2104
$ws.Range("I137").Value = 2648.75
$ws.Range("K137").Value = 7946.25
$ws.Range("M137").Value = -5396.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4237.706
$ws.Range("I2").Value = 3407.75
$ws.Range("K2").Value = 3407.75
$ws.Range("M2").Value = -3294.75
$ws.Range("H44").Value = 14609.8
$ws.Range("J44").Value = 15512.25
$ws.Range("L44").Value = 15512.25
$ws.Range("N44").Value = -16488.25
$ws.Range("H45").Value = 4115.6875
$ws.Range("I45").Value = 1722.7778
$ws.Range("K45").Value = 1722.7778
$ws.Range("M45").Value = -1345.7778
$ws.Range("H55").Value = 24017.666
$ws.Range("J55").Value = 23526.5
$ws.Range("L55").Value = 23526.5
$ws.Range("N55").Value = -24156.5
$ws.Range("H74").Value = 11095.483
$ws.Range("J74").Value = 5342
$ws.Range("L74").Value = 5342
$ws.Range("N74").Value = -7090
$ws.Range("H77").Value = 11095.483
$ws.Range("J77").Value = 5342
$ws.Range("L77").Value = 26710
$ws.Range("N77").Value = -35446
$ws.Range("H110").Value = 839002
$ws.Range("I110").Value = 1004799.8
$ws.Range("K110").Value = 1004799.8
$ws.Range("M110").Value = -1002754.8
$ws.Range("H116").Value = 4237.706
$ws.Range("I116").Value = 3407.75
$ws.Range("K116").Value = 3407.75
$ws.Range("M116").Value = -1113.75
$ws.Range("H122").Value = 3486.2903
$ws.Range("I122").Value = 3034.3076
$ws.Range("K122").Value = 9102.9228
$ws.Range("M122").Value = -6652.9228
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4237.706
$ws.Range("I3").Value = 3407.75
$ws.Range("K3").Value = 3407.75
$ws.Range("M3").Value = -3293.75
$ws.Range("H80").Value = 2139.6875
$ws.Range("J80").Value = 1805
$ws.Range("L80").Value = 1805
$ws.Range("N80").Value = -3801
$ws.Range("H83").Value = 2139.6875
$ws.Range("J83").Value = 1805
$ws.Range("L83").Value = 9025
$ws.Range("N83").Value = -19009
$ws.Range("H105").Value = 2041.7778
$ws.Range("I105").Value = 1235.6
$ws.Range("J105").Value = 2351.8462
$ws.Range("K105").Value = 1235.6
$ws.Range("L105").Value = 2351.8462
$ws.Range("M105").Value = 511.4000000000001
$ws.Range("N105").Value = -5845.8462

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3168.889
$ws.Range("I31").Value = 2020.6364
$ws.Range("J31").Value = 4973.2856
$ws.Range("K31").Value = 2020.6364
$ws.Range("L31").Value = 4973.2856
$ws.Range("M31").Value = -1725.6364
$ws.Range("N31").Value = -5563.2856
$ws.Range("H34").Value = 3168.889
$ws.Range("I34").Value = 2020.6364
$ws.Range("J34").Value = 4973.2856
$ws.Range("K34").Value = 2020.6364
$ws.Range("L34").Value = 4973.2856
$ws.Range("M34").Value = -1818.6364
$ws.Range("N34").Value = -5377.2856
$ws.Range("H58").Value = 230587.95
$ws.Range("I58").Value = 437415.4
$ws.Range("K58").Value = 437415.4
$ws.Range("M58").Value = -437212.4
$ws.Range("H134").Value = 3104.4285
$ws.Range("I134").Value = 2403.0344
$ws.Range("K134").Value = 7209.1032
$ws.Range("M134").Value = -4674.1032
$ws.Range("H136").Value = 230587.95
$ws.Range("I136").Value = 437415.4
$ws.Range("K136").Value = 1312246.2
$ws.Range("M136").Value = -1309696.2

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 607
$ws.Range("J75").Value = 554.5
$ws.Range("L75").Value = 1663.5
$ws.Range("N75").Value = -3659.5
$ws.Range("H78").Value = 607
$ws.Range("J78").Value = 554.5
$ws.Range("L78").Value = 4990.5
$ws.Range("N78").Value = -14974.5
$ws.Range("H121").Value = 2000660
$ws.Range("I121").Value = 500
$ws.Range("K121").Value = 1500
$ws.Range("M121").Value = -190
$ws.Range("H131").Value = 3092.5527
$ws.Range("J131").Value = 4110.0835
$ws.Range("L131").Value = 12330.2505
$ws.Range("N131").Value = -22410.2505
$ws.Range("H132").Value = 2057.2
$ws.Range("I132").Value = 779.55554
$ws.Range("J132").Value = 5342.5713
$ws.Range("K132").Value = 7015.99986
$ws.Range("L132").Value = 48083.14169999999
$ws.Range("M132").Value = -4485.99986
$ws.Range("N132").Value = -53143.14169999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 40000
$ws.Range("J20").Value = 40000
$ws.Range("L20").Value = 40000
$ws.Range("N20").Value = -40490
$ws.Range("H51").Value = 59441.668
$ws.Range("J51").Value = 59441.668
$ws.Range("L51").Value = 59441.668
$ws.Range("N51").Value = -60459.668
$ws.Range("H113").Value = 1669980.1
$ws.Range("J113").Value = 3967.5
$ws.Range("L113").Value = 3967.5
$ws.Range("N113").Value = -8307.5
$ws.Range("H140").Value = 49999.832
$ws.Range("J140").Value = 49999.832
$ws.Range("L140").Value = 49999.832
$ws.Range("N140").Value = -60359.832

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 532.2059
$ws.Range("I55").Value = 321.68
$ws.Range("K55").Value = 321.68
$ws.Range("M55").Value = -148.68
$ws.Range("H61").Value = 4872.6
$ws.Range("I61").Value = 3557
$ws.Range("K61").Value = 3557
$ws.Range("M61").Value = -3355
$ws.Range("H108").Value = 34999.5
$ws.Range("J108").Value = 34999.5
$ws.Range("L108").Value = 34999.5
$ws.Range("N108").Value = -42679.5
$ws.Range("H113").Value = 4872.6
$ws.Range("I113").Value = 3557
$ws.Range("K113").Value = 3557
$ws.Range("M113").Value = -1387
$ws.Range("H122").Value = 442257.47
$ws.Range("I122").Value = 316795.75
$ws.Range("J122").Value = 693180.9
$ws.Range("K122").Value = 950387.25
$ws.Range("L122").Value = 2079542.7
$ws.Range("M122").Value = -947937.25
$ws.Range("N122").Value = -2084442.7
$ws.Range("H136").Value = 4905.8
$ws.Range("I136").Value = 3071.3333
$ws.Range("K136").Value = 9213.999899999999
$ws.Range("M136").Value = -6663.999899999999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 79198.08
$ws.Range("I107").Value = 85589.586
$ws.Range("J107").Value = 2500
$ws.Range("K107").Value = 256768.758
$ws.Range("L107").Value = 7500
$ws.Range("M107").Value = -254848.758
$ws.Range("N107").Value = -11340
$ws.Range("H113").Value = 1990.9286
$ws.Range("I113").Value = 2110.875
$ws.Range("J113").Value = 1831
$ws.Range("K113").Value = 6332.625
$ws.Range("L113").Value = 5493
$ws.Range("M113").Value = -4162.625
$ws.Range("N113").Value = -9833
$ws.Range("H122").Value = 21279786
$ws.Range("I122").Value = 32260364
$ws.Range("K122").Value = 96781092
$ws.Range("M122").Value = -96778642
$ws.Range("H136").Value = 360736.56
$ws.Range("I136").Value = 478792.8
$ws.Range("K136").Value = 1436378.4
$ws.Range("M136").Value = -1433828.4

